$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 349, pushing existing rows 349:415 down to 350:416
$ws.Rows.Item(349).Insert()

# Populate the newly inserted row 349 with the new weekly price record
$ws.Cells.Item(349, 1).Value  = 3
$ws.Cells.Item(349, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(349, 3).Value  = "Coquimbo"
$ws.Cells.Item(349, 4).Value  = 45258
$ws.Cells.Item(349, 5).Value  = 5
$ws.Cells.Item(349, 6).Value  = "Fruta"
$ws.Cells.Item(349, 7).Value  = 100101
$ws.Cells.Item(349, 8).Value  = "Berries"
$ws.Cells.Item(349, 9).Value  = 100101001
$ws.Cells.Item(349, 10).Value = "Arándano (blue)"
$ws.Cells.Item(349, 11).Value = "Sin especificar"
$ws.Cells.Item(349, 12).Value = "Primera"
$ws.Cells.Item(349, 13).Value = 45
$ws.Cells.Item(349, 14).Value = 5600
$ws.Cells.Item(349, 15).Value = 5600
$ws.Cells.Item(349, 16).Value = 5600
$ws.Cells.Item(349, 17).Value = "$/bandeja 2 kilos"
$ws.Cells.Item(349, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(349, 19).Value = 2800
$ws.Cells.Item(349, 20).Value = 2
